# Atualizado por script em 11-11-2023 14:45
#
# Corrects two home/away row-ordering swaps (rows 23/24, 60/61, 74/75,
# 77/78 each had their match data stored against the wrong row index)
# and appends two newly scraped fixtures (rows 84 and 85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# swap rows 23 and 24
$ws.Range("F23").Value = 'Oliveirense'
$ws.Range("F24").Value = 'FC Porto B'
$ws.Range("G23").Value = 3
$ws.Range("G24").Value = 2
$ws.Range("H23").Value = 'Penafiel'
$ws.Range("H24").Value = 'Leiria'
$ws.Range("I23").Value = 1
$ws.Range("I24").Value = 1
$ws.Range("J23").Value = 2.28
$ws.Range("J24").Value = 2.17
$ws.Range("K23").Value = '23/08/2023 10:41'
$ws.Range("K24").Value = '23/08/2023 10:41'
$ws.Range("L23").Value = 2.42
$ws.Range("L24").Value = 2.7
$ws.Range("M23").Value = '27/08/2023 11:40'
$ws.Range("M24").Value = '27/08/2023 11:51'
$ws.Range("N23").Value = 3.34
$ws.Range("N24").Value = 3.44
$ws.Range("O23").Value = '23/08/2023 10:41'
$ws.Range("O24").Value = '23/08/2023 10:41'
$ws.Range("P23").Value = 3.46
$ws.Range("P24").Value = 3.37
$ws.Range("Q23").Value = '27/08/2023 10:00'
$ws.Range("Q24").Value = '27/08/2023 11:51'
$ws.Range("R23").Value = 3.29
$ws.Range("R24").Value = 3.43
$ws.Range("S23").Value = '23/08/2023 10:41'
$ws.Range("S24").Value = '23/08/2023 10:41'
$ws.Range("T23").Value = 3.03
$ws.Range("T24").Value = 2.74
$ws.Range("U23").Value = '27/08/2023 11:40'
$ws.Range("U24").Value = '27/08/2023 11:51'
$ws.Range("V23").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/oliveirense-penafiel/h4QIwcVj/'
$ws.Range("V24").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-leiria/z1LGbb7G/'

# swap rows 60 and 61
$ws.Range("F60").Value = 'Mafra'
$ws.Range("F61").Value = 'FC Porto B'
$ws.Range("G60").Value = 3
$ws.Range("G61").Value = 2
$ws.Range("H60").Value = 'Leiria'
$ws.Range("H61").Value = 'Torreense'
$ws.Range("I60").Value = 0
$ws.Range("I61").Value = 2
$ws.Range("J60").Value = 2.28
$ws.Range("J61").Value = 2.36
$ws.Range("K60").Value = '04/10/2023 14:42'
$ws.Range("K61").Value = '04/10/2023 14:42'
$ws.Range("L60").Value = 2.86
$ws.Range("L61").Value = 2.73
$ws.Range("M60").Value = '08/10/2023 11:59'
$ws.Range("M61").Value = '08/10/2023 11:51'
$ws.Range("N60").Value = 3.48
$ws.Range("N61").Value = 3.38
$ws.Range("O60").Value = '04/10/2023 14:42'
$ws.Range("O61").Value = '04/10/2023 14:42'
$ws.Range("P60").Value = 3.41
$ws.Range("P61").Value = 3.34
$ws.Range("Q60").Value = '08/10/2023 11:59'
$ws.Range("Q61").Value = '08/10/2023 11:53'
$ws.Range("R60").Value = 2.98
$ws.Range("R61").Value = 2.93
$ws.Range("S60").Value = '04/10/2023 14:42'
$ws.Range("S61").Value = '04/10/2023 14:42'
$ws.Range("T60").Value = 2.57
$ws.Range("T61").Value = 2.73
$ws.Range("U60").Value = '08/10/2023 11:59'
$ws.Range("U61").Value = '08/10/2023 11:51'
$ws.Range("V60").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-leiria/hbqtOZiP/'
$ws.Range("V61").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-torreense/IL4w4E0g/'

# swap rows 74 and 75
$ws.Range("F74").Value = 'Mafra'
$ws.Range("F75").Value = 'FC Porto B'
$ws.Range("G74").Value = 0
$ws.Range("G75").Value = 2
$ws.Range("H74").Value = 'Leixoes'
$ws.Range("H75").Value = 'Feirense'
$ws.Range("I74").Value = 1
$ws.Range("I75").Value = 0
$ws.Range("J74").Value = 1.88
$ws.Range("J75").Value = 1.98
$ws.Range("K74").Value = '01/11/2023 16:12'
$ws.Range("K75").Value = '01/11/2023 16:12'
$ws.Range("L74").Value = 1.93
$ws.Range("L75").Value = 1.84
$ws.Range("M74").Value = '04/11/2023 11:48'
$ws.Range("M75").Value = '04/11/2023 11:59'
$ws.Range("N74").Value = 3.73
$ws.Range("N75").Value = 3.57
$ws.Range("O74").Value = '01/11/2023 16:12'
$ws.Range("O75").Value = '01/11/2023 16:12'
$ws.Range("P74").Value = 3.55
$ws.Range("P75").Value = 3.78
$ws.Range("Q74").Value = '04/11/2023 11:51'
$ws.Range("Q75").Value = '04/11/2023 11:59'
$ws.Range("R74").Value = 3.8
$ws.Range("R75").Value = 3.87
$ws.Range("S74").Value = '01/11/2023 16:12'
$ws.Range("S75").Value = '01/11/2023 16:12'
$ws.Range("T74").Value = 4.21
$ws.Range("T75").Value = 4.41
$ws.Range("U74").Value = '04/11/2023 11:51'
$ws.Range("U75").Value = '04/11/2023 11:58'
$ws.Range("V74").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-leixoes/YiBBPnTT/'
$ws.Range("V75").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-feirense/jTL6QSDN/'

# swap rows 77 and 78
$ws.Range("F77").Value = 'Benfica B'
$ws.Range("F78").Value = 'Nacional'
$ws.Range("G77").Value = 1
$ws.Range("G78").Value = 1
$ws.Range("H77").Value = 'Penafiel'
$ws.Range("H78").Value = 'Santa Clara'
$ws.Range("I77").Value = 0
$ws.Range("I78").Value = 1
$ws.Range("J77").Value = 1.91
$ws.Range("J78").Value = 2.98
$ws.Range("K77").Value = '29/10/2023 16:42'
$ws.Range("K78").Value = '01/11/2023 16:12'
$ws.Range("L77").Value = 2.26
$ws.Range("L78").Value = 2.81
$ws.Range("M77").Value = '04/11/2023 18:53'
$ws.Range("M78").Value = '04/11/2023 18:58'
$ws.Range("N77").Value = 3.67
$ws.Range("N78").Value = 3.27
$ws.Range("O77").Value = '29/10/2023 16:42'
$ws.Range("O78").Value = '01/11/2023 16:12'
$ws.Range("P77").Value = 3.59
$ws.Range("P78").Value = 3.23
$ws.Range("Q77").Value = '04/11/2023 18:53'
$ws.Range("Q78").Value = '04/11/2023 18:52'
$ws.Range("R77").Value = 4.01
$ws.Range("R78").Value = 2.39
$ws.Range("S77").Value = '29/10/2023 16:42'
$ws.Range("S78").Value = '01/11/2023 16:12'
$ws.Range("T77").Value = 3.2
$ws.Range("T78").Value = 2.72
$ws.Range("U77").Value = '04/11/2023 18:53'
$ws.Range("U78").Value = '04/11/2023 18:52'
$ws.Range("V77").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/benfica-penafiel/xjmbUAEb/'
$ws.Range("V78").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-santa-clara/xQH2R8bH/'
# Append new rows 84 and 85, copying number-format/style from row 83 first
$ws.Range("A83:V83").Copy()
$ws.Range("A84:V85").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 84
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = 'portugal'
$ws.Range("C84").Value = 'liga-portugal-2'
$ws.Range("D84").Value = '2023-2024'
$ws.Range("E84").Value = 45241.5
$ws.Range("F84").Value = 'Nacional'
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 'Leiria'
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 2.52
$ws.Range("K84").Value = '08/11/2023 06:12'
$ws.Range("L84").Value = 2.32
$ws.Range("M84").Value = '11/11/2023 11:58'
$ws.Range("N84").Value = 3.42
$ws.Range("O84").Value = '08/11/2023 06:12'
$ws.Range("P84").Value = 3.59
$ws.Range("Q84").Value = '11/11/2023 11:52'
$ws.Range("R84").Value = 2.83
$ws.Range("S84").Value = '08/11/2023 06:12'
$ws.Range("T84").Value = 3.09
$ws.Range("U84").Value = '11/11/2023 11:58'
$ws.Range("V84").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-leiria/Qozkw6cT/'

# row 85
$ws.Range("A85").Value = 84
$ws.Range("B85").Value = 'portugal'
$ws.Range("C85").Value = 'liga-portugal-2'
$ws.Range("D85").Value = '2023-2024'
$ws.Range("E85").Value = 45241.625
$ws.Range("F85").Value = 'Penafiel'
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 'FC Porto B'
$ws.Range("I85").Value = 2
$ws.Range("J85").Value = 2.66
$ws.Range("K85").Value = '08/11/2023 06:12'
$ws.Range("L85").Value = 3.48
$ws.Range("M85").Value = '11/11/2023 14:54'
$ws.Range("N85").Value = 3.32
$ws.Range("O85").Value = '08/11/2023 06:12'
$ws.Range("P85").Value = 3.33
$ws.Range("Q85").Value = '11/11/2023 14:54'
$ws.Range("R85").Value = 2.74
$ws.Range("S85").Value = '08/11/2023 06:12'
$ws.Range("T85").Value = 2.24
$ws.Range("U85").Value = '11/11/2023 14:54'
$ws.Range("V85").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal-2/penafiel-fc-porto/27kPrATi/'

Write-Host "Done. UsedRange is now:" $ws.UsedRange.Address()
